$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header block (row 4): project name / version / date ---
$ws.Range("C4").Value = "Développement mobile"
$ws.Range("E4").Value = "version: 1"
$ws.Range("G4").Value = "date: 22/03/2023"

# --- Student names (row 6) ---
$ws.Range("D6").Value = "Gauzi"
$ws.Range("E6").Value = "Paillard"
$ws.Range("F6").ClearContents()
$ws.Range("G6").ClearContents()
$ws.Range("H6").ClearContents()

# --- Task rows ---
# Row 7
$ws.Range("C7").Value = "Déploiment du serveur ASP.NET"
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 1
$ws.Range("F7").ClearContents()

# Row 8
$ws.Range("C8").Value = "Rédaction du rapport"
$ws.Range("D8").Value = 0.4
$ws.Range("E8").Value = 0.6

# Row 9
$ws.Range("C9").Value = "Intialisation du projet "
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 0

# Row 10
$ws.Range("C10").ClearContents()

# --- Selection ---
$ws.Range("E11").Select() | Out-Null

Write-Output "done"
